$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (2014/12) ----
$ws.Range("D2").Value = 13050
$ws.Range("E2").Value = 578
$ws.Range("F2").Value = 578
$ws.Range("G2").Value = 430
$ws.Range("H2").Value = 317
$ws.Range("I2").Value = 217
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 12125
$ws.Range("L2").Value = 7820
$ws.Range("M2").Value = 4304
$ws.Range("N2").Value = 2844
$ws.Range("O2").Value = 1460
$ws.Range("P2").Value = 250
$ws.Range("Q2").Value = 752
$ws.Range("R2").Value = -437
$ws.Range("S2").Value = -694
$ws.Range("T2").Value = 216
$ws.Range("U2").Value = 535
$ws.Range("V2").Value = 5124
$ws.Range("W2").Value = 4.43
$ws.Range("X2").Value = 2.43
$ws.Range("Y2").Value = 8.09
$ws.Range("Z2").Value = 2.38
$ws.Range("AA2").Value = 181.69
$ws.Range("AB2").Value = 1042.56
$ws.Range("AC2").Value = 4348
$ws.Range("AD2").Value = 9.15
$ws.Range("AE2").Value = 57522
$ws.Range("AF2").Value = 0.6899999999999999
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 5000000

# ---- Row 3 (2015/12) ----
$ws.Range("D3").Value = 9052
$ws.Range("E3").Value = 341
$ws.Range("F3").Value = 341
$ws.Range("G3").Value = 363
$ws.Range("H3").Value = 244
$ws.Range("I3").Value = 191
$ws.Range("J3").Value = 53
$ws.Range("K3").Value = 9031
$ws.Range("L3").Value = 5253
$ws.Range("M3").Value = 3778
$ws.Range("N3").Value = 3096
$ws.Range("O3").Value = 682
$ws.Range("P3").Value = 250
$ws.Range("Q3").Value = 359
$ws.Range("R3").Value = 84
$ws.Range("S3").Value = -471
$ws.Range("T3").Value = 151
$ws.Range("U3").Value = 208
$ws.Range("V3").Value = 3264
$ws.Range("W3").Value = 3.77
$ws.Range("X3").Value = 2.69
$ws.Range("Y3").Value = 6.43
$ws.Range("Z3").Value = 2.31
$ws.Range("AA3").Value = 139.06
$ws.Range("AB3").Value = 1100.03
$ws.Range("AC3").Value = 3821
$ws.Range("AD3").Value = 17.69
$ws.Range("AE3").Value = 62604
$ws.Range("AF3").Value = 1.08
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 5000000

# ---- Row 4 (2016/12) ----
$ws.Range("D4").Value = 7021
$ws.Range("E4").Value = 492
$ws.Range("F4").Value = 492
$ws.Range("G4").Value = 427
$ws.Range("H4").Value = 301
$ws.Range("I4").Value = 264
$ws.Range("J4").Value = 36
$ws.Range("K4").Value = 10107
$ws.Range("L4").Value = 6030
$ws.Range("M4").Value = 4077
$ws.Range("N4").Value = 3361
$ws.Range("O4").Value = 716
$ws.Range("P4").Value = 250
$ws.Range("Q4").Value = 383
$ws.Range("R4").Value = -1073
$ws.Range("S4").Value = 646
$ws.Range("T4").Value = 381
$ws.Range("U4").Value = 3
$ws.Range("V4").Value = 3910
$ws.Range("W4").Value = 7.01
$ws.Range("X4").Value = 4.28
$ws.Range("Y4").Value = 8.18
$ws.Range("Z4").Value = 3.14
$ws.Range("AA4").Value = 147.91
$ws.Range("AB4").Value = 1202.28
$ws.Range("AC4").Value = 5282
$ws.Range("AD4").Value = 12.1
$ws.Range("AE4").Value = 67969
$ws.Range("AF4").Value = 0.9399999999999999
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 0.23
$ws.Range("AI4").Value = 2.81
$ws.Range("AJ4").Value = 5000000

# ---- Row 5 (2017/12) ----
$ws.Range("D5").Value = 8160
$ws.Range("E5").Value = 591
$ws.Range("F5").Value = 591
$ws.Range("G5").Value = 742
$ws.Range("H5").Value = 541
$ws.Range("I5").Value = 437
$ws.Range("J5").Value = 104
$ws.Range("K5").Value = 10573
$ws.Range("L5").Value = 5917
$ws.Range("M5").Value = 4656
$ws.Range("N5").Value = 3823
$ws.Range("O5").Value = 833
$ws.Range("P5").Value = 250
$ws.Range("Q5").Value = 664
$ws.Range("R5").Value = -165
$ws.Range("S5").Value = -460
$ws.Range("T5").Value = 211
$ws.Range("U5").Value = 453
$ws.Range("V5").Value = 3583
$ws.Range("W5").Value = 7.25
$ws.Range("X5").Value = 6.63
$ws.Range("Y5").Value = 12.18
$ws.Range("Z5").Value = 5.24
$ws.Range("AA5").Value = 127.07
$ws.Range("AB5").Value = 1375.97
$ws.Range("AC5").Value = 8750
$ws.Range("AD5").Value = 7.81
$ws.Range("AE5").Value = 77316
$ws.Range("AF5").Value = 0.88
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 0.29
$ws.Range("AI5").Value = 2.26
$ws.Range("AJ5").Value = 5000000

# ---- Row 6 (2018/12) ----
$ws.Range("D6").Value = 7821
$ws.Range("E6").Value = 543
$ws.Range("F6").Value = 543
$ws.Range("G6").Value = 563
$ws.Range("H6").Value = 412
$ws.Range("I6").Value = 356
$ws.Range("K6").Value = 11760
$ws.Range("L6").Value = 6662
$ws.Range("M6").Value = 5098
$ws.Range("N6").Value = 4176
$ws.Range("P6").Value = 250
$ws.Range("Q6").Value = 164
$ws.Range("R6").Value = -351
$ws.Range("S6").Value = 468
$ws.Range("T6").Value = 488
$ws.Range("U6").Value = -324
$ws.Range("V6").Value = 4158
$ws.Range("W6").Value = 6.94
$ws.Range("X6").Value = 5.27
$ws.Range("Y6").Value = 8.9
$ws.Range("Z6").Value = 3.69
$ws.Range("AA6").Value = 130.67
$ws.Range("AB6").Value = 1511.88
$ws.Range("AC6").Value = 7123
$ws.Range("AD6").Value = 7
$ws.Range("AE6").Value = 84456
$ws.Range("AF6").Value = 0.59
$ws.Range("AI6").Value = 3.47
$ws.Range("AJ6").Value = 5000000

# AG6 / AH6 no longer exist in the updated data (2019/12(E) row lost its dividend figures)
$ws.Range("AG6:AH6").ClearContents()

# ---- Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) ----
# These forward estimate rows are cleared entirely except for the A/B/C label columns
$ws.Range("D7:AJ9").ClearContents()
